$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$colB = @(13.8857373893269, 13.41996694170678, 13.13171107950289, 13.01390084365406, 12.99432473054772, 13.13012332911635, 13.72574269336042, 14.86655597787522, 15.67705597860293, 16.03777742432541, 16.17308322880798, 16.14400226942056, 16.04893576753983, 15.99053258736239, 15.65330872314058, 15.44428133538532, 15.32331083028994, 15.28222926589065, 15.46661070161354, 16.07689522513528, 16.46816753970792, 16.26007554401473, 15.45651806922899, 14.562058765201)
$colC = @(12.70541594692189, 12.03972334951779, 11.61582460479961, 11.43950254710403, 11.41001575853982, 11.61346082259347, 12.47915941559864, 14.04813205869586, 15.11342594154173, 15.57769849600122, 15.75049549574828, 15.71341580248056, 15.59197544009676, 15.51719499957836, 15.08266745081519, 14.81082169109396, 14.65255324996099, 14.59864123822366, 14.83995850792939, 15.62772781054443, 16.12498934118858, 15.86122608644919, 14.82679191553898, 13.63849473468293)
$colE = @(10.01737787601544, 9.899375052573539, 9.829073681116677, 9.800993428307571, 9.79636579220265, 9.828692646322464, 9.976260154150971, 10.28148489177095, 10.51367214559231, 10.62064572848342, 10.66131598117923, 10.65255023996512, 10.62398864007859, 10.60651392814549, 10.50670556303444, 10.44579908717324, 10.41089635035693, 10.39910204116359, 10.4522695644471, 10.63237374479751, 10.7510096947389, 10.68761728588628, 10.44934390658804, 10.19739200742938)
$colF = @(16.86991607391245, 15.89584955866815, 15.26997757108489, 15.008197319934, 14.96433081551589, 15.26647399323726, 16.53996406344765, 19.0027458068253, 20.67494806633232, 21.3917225636224, 21.65686569030329, 21.60004134736742, 21.4136618050453, 21.29868154950795, 20.62722412089977, 20.20408069617459, 19.95656407809808, 19.87204792380562, 20.2495528364879, 21.46857628470567, 22.22866616901555, 21.82633154475864, 20.22900810905294, 18.34778573295697)
$colG = @(3.646307196277861, 3.650925257400409, 3.653899826192809, 3.655147117978964, 3.655356356415625, 3.653916505145128, 3.647870740154809, 3.637110939138823, 3.629863342777296, 3.626706764953191, 3.625531463040794, 3.625783697280123, 3.626609671746018, 3.627118208056464, 3.630072443701285, 3.631920613463786, 3.632996855740166, 3.633363529038573, 3.631722505188332, 3.626366520829478, 3.62298272436533, 3.624778100562143, 3.631812027239994, 3.639905506164345)
$colI = @(29.41626950358807, 28.87340025757983, 28.5396989231449, 28.40375149824199, 28.38118377560937, 28.53786515554759, 29.22924628024607, 30.57626426675839, 31.55301568879818, 31.99303838699155, 32.15892315375176, 32.12323181168701, 32.00670138717946, 31.93522291747164, 31.52416186893206, 31.27079694341799, 31.12466981599801, 31.07512910970709, 31.2978102064381, 32.04095029473778, 32.52225581782003, 32.26581386055894, 31.2855989468867, 30.21359276988478)
$colM = @(16.09831625193257, 15.84666363530523, 15.69591011986492, 15.63550409830809, 15.62553809925424, 15.6950911980338, 16.01081491098055, 16.656150422034, 17.14132046254653, 17.3633974006449, 17.44760483588757, 17.42946562173176, 17.37032343808485, 17.33410920069087, 17.12682743961227, 16.99995635506587, 16.92711745082659, 16.90248110255404, 17.01344875762033, 17.38769256695668, 17.63289600381555, 17.50199771110729, 17.0073485214142, 16.4792746098998)

$cols = @{ "B" = $colB; "C" = $colC; "E" = $colE; "F" = $colF; "G" = $colG; "I" = $colI; "M" = $colM }

foreach ($colLetter in $cols.Keys) {
    $values = $cols[$colLetter]
    $data = New-Object 'object[,]' $values.Count, 1
    for ($i = 0; $i -lt $values.Count; $i++) {
        $data[$i, 0] = $values[$i]
    }
    $startRow = 2
    $endRow = $startRow + $values.Count - 1
    $range = $ws.Range("$colLetter$startRow`:$colLetter$endRow")
    $range.Value = $data
}
